# Rename the existing sheet "Sayfa1" -> "kodlar" and add a new sheet
# "adminler" after it, populated with admin/login style data, matching
# the target commit.

$wb = $excel.ActiveWorkbook

# --- Rename the original sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "kodlar"

# --- Add the new sheet right after "kodlar" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "adminler"

# Copy the formatting used on "kodlar" onto the new sheet so the new
# cells share the same cell style as the rest of the workbook instead of
# getting a brand new style definition.
$ws1.Range("A1:B3").Copy()
$ws2.Range("A1:B3").PasteSpecial(-4122) # xlPasteFormats

# --- Fill in the data for "adminler" ---
$ws2.Range("A1").Value = "kodlar"
$ws2.Range("B1").Value = "isim"

$ws2.Range("A2").Value = 4444
$ws2.Range("B2").Value = "Mert"

$ws2.Range("A3").Value = 3333
$ws2.Range("B3").Value = "hasan"
